# Weekly update: add this week's "Brócoli" price rows at the top of the
# price-history block (rows 1136-1164), shifting the existing rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block (pushes old
# row 1136.. down to 1138..1166, and auto-extends the date style/format
# from the surrounding rows onto the new rows' cells).
$ws.Rows("1136:1137").Insert()

# New row 1136: "Primera" quality entry for the new week (2023-11-09)
$ws.Range("A1136").Value = 9
$ws.Range("B1136").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1136").Value = "Metropolitana"
$ws.Range("D1136").Value = 45239
$ws.Range("E1136").Value = 13
$ws.Range("F1136").Value = 100112023
$ws.Range("G1136").Value = "Brócoli"
$ws.Range("H1136").Value = "Sin especificar"
$ws.Range("I1136").Value = "Primera"
$ws.Range("J1136").Value = 3400
$ws.Range("K1136").Value = 800
$ws.Range("L1136").Value = 900
$ws.Range("M1136").Value = 850
$ws.Range("N1136").Value = "$/unidad"
$ws.Range("O1136").Value = "Región Metropolitana"
$ws.Range("P1136").Value = 850
$ws.Range("Q1136").Value = 1
$ws.Range("R1136").Value = "Hortaliza"

# New row 1137: "Segunda" quality entry for the same new week
$ws.Range("A1137").Value = 9
$ws.Range("B1137").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1137").Value = "Metropolitana"
$ws.Range("D1137").Value = 45239
$ws.Range("E1137").Value = 13
$ws.Range("F1137").Value = 100112023
$ws.Range("G1137").Value = "Brócoli"
$ws.Range("H1137").Value = "Sin especificar"
$ws.Range("I1137").Value = "Segunda"
$ws.Range("J1137").Value = 1690
$ws.Range("K1137").Value = 800
$ws.Range("L1137").Value = 800
$ws.Range("M1137").Value = 800
$ws.Range("N1137").Value = "$/unidad"
$ws.Range("O1137").Value = "Región Metropolitana"
$ws.Range("P1137").Value = 800
$ws.Range("Q1137").Value = 1
$ws.Range("R1137").Value = "Hortaliza"
